$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert new BOM line for feedthrough caps (C586-609) above current row 106 ---
$ws.Rows.Item(106).Insert()
$ws.Cells.Item(106,1).Value = 105
$ws.Cells.Item(106,2).Value = 24
$ws.Cells.Item(106,3).Value = "C586-609"
$ws.Cells.Item(106,4).Value = "NFM18PC105R0J3D"
$ws.Cells.Item(106,5).Value = "NFM18PC105R0J3D"
$ws.Cells.Item(106,6).Value = "Murata"
$ws.Cells.Item(106,7).Value = "NFM18PC105R0J3D"
$ws.Cells.Item(106,8).Value = "EMIFIL FILTER  NFM18P Murata NFM18C Series, Signal Filter, 6.3 V dc, 2A 0603 SMD 1.6 x 0.8 x 0.8mm"

# --- Insert new BOM line for 2mm test points (TP128-175) above current row 166 ---
# (this row shifted down by 1 because of the insert above)
$ws.Rows.Item(166).Insert()
$ws.Cells.Item(166,1).Value = 165
$ws.Cells.Item(166,2).Value = 48
$ws.Cells.Item(166,3).Value = "TP128-175"
$ws.Cells.Item(166,4).Value = "TESTPOINT_2MM"
$ws.Cells.Item(166,5).Value = "TP_2MM"
$ws.Cells.Item(166,6).Value = "DNI"
$ws.Cells.Item(166,7).Value = "DNI"
$ws.Cells.Item(166,8).Value = "Test Point"

# --- Renumber the "Line #" column (A) for every data row so it stays row-1 ---
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r,1).Value = $r - 1
}

# --- Update selection / scroll position to match the post-edit view ---
$ws.Application.Goto($ws.Range("A101"), $false)
$ws.Range("G106").Select()
